$d = $word.ActiveDocument

# 1. Update the letter date: "May 6, 2019" -> "May 23, 2019"
$ok1 = $d.Content.Find.Execute("May 6, 2019", $true, $false, $false, $false, $false,
                                $true, 1, $false, "May 23, 2019", 2)

# 2. Re-word the manuscript-title reference:
#    ' tag loss model" as an expository paper' -> ' model with tag loss"'
$ok2 = $d.Content.Find.Execute(" tag loss model” as an expository paper",
                                $true, $false, $false, $false, $false,
                                $true, 1, $false,
                                " model with tag loss”",
                                2)

# 3. Add the missing period after "involving elephant seals"
$ok3 = $d.Content.Find.Execute("involving elephant seals  Recycled",
                                $true, $false, $false, $false, $false,
                                $true, 1, $false,
                                "involving elephant seals.  Recycled", 2)

# 4. Punctuation fix: "tagged thus this issue" -> "tagged; thus, this issue",
#    plus the extra space that now precedes "However" (double space after the
#    previous sentence, matching the rest of the letter's style).
$ok4 = $d.Content.Find.Execute(
    "individuals. However, in the elephant seal data, individuals were branded as well as tagged thus this issue could be explored and we offer",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "individuals.  However, in the elephant seal data, individuals were branded as well as tagged; thus, this issue could be explored and we offer",
    2)

Write-Output ("Replacements ok: " + $ok1 + " " + $ok2 + " " + $ok3 + " " + $ok4)

# 5. Move the "_GoBack" bookmark from the end of the signature paragraph to the
#    final (empty) paragraph of the document.
$bookmarks = $d.Bookmarks
if ($bookmarks.Exists("_GoBack")) {
    $bookmarks.Item("_GoBack").Delete()
}
$lastPara = $d.Paragraphs.Last
$bookmarks.Add("_GoBack", $lastPara.Range)
Write-Output ("_GoBack now exists: " + $bookmarks.Exists("_GoBack"))
